$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refactor state labels by group -----------------------------------
# Capitalize the "(away)"/"(stay)" qualifiers consistently, and rename
# "Breach" -> "Breached". Every cell that displayed the old label needs
# updating (header row + the A-column state names + the body cells that
# reference those states as sources/targets).

# Header row: "Timer Ticks" / "Timer Runs Out" swap columns I and J
$ws.Range("I1").Value = "Timer Runs Out"
$ws.Range("J1").Value = "Timer Ticks"

# Arming (away)/(stay) -> Arming (Away)/(Stay)
$ws.Range("B2").Value = "Arming (Away)"
$ws.Range("C2").Value = "Arming (Stay)"
$ws.Range("A3").Value = "Arming (Away)"
$ws.Range("A4").Value = "Arming (Stay)"

# Confirm (away)/(stay) -> Confirm (Away)/(Stay); these moved from column J
# (Confirm became a self-contained row entry) into column I on the rows
# for the Arming states.
$ws.Range("I3").Value = "Confirm (Away)"
$ws.Range("J3").Value = "~"
$ws.Range("I4").Value = "Confirm (Stay)"
$ws.Range("J4").Value = "~"
$ws.Range("A5").Value = "Confirm (Away)"
$ws.Range("A6").Value = "Confirm (Stay)"

# Armed (away)/(stay) -> Armed (Away)/(Stay)
$ws.Range("E5").Value = "Armed (Away)"
$ws.Range("E6").Value = "Armed (Stay)"
$ws.Range("A7").Value = "Armed (Away)"
$ws.Range("A8").Value = "Armed (Stay)"

# Disarm (away)/(stay) -> Disarm (Away)/(Stay)
$ws.Range("D7").Value = "Disarm (Away)"
$ws.Range("D8").Value = "Disarm (Stay)"

# Breach -> Breached
$ws.Range("F8").Value = "Breached"

# --- Reorder the remaining states by group (Disarm, Warning, Breached) -
$ws.Range("A9").Value = "Disarm (Away)"
$ws.Range("I9").Value = "Armed (Away)"
$ws.Range("J9").Value = "~"

$ws.Range("A10").Value = "Disarm (Stay)"
$ws.Range("I10").Value = "Armed (Stay)"

$ws.Range("A11").Value = "Warning"
$ws.Range("I11").Value = "Breached"
$ws.Range("J11").Value = "~"

$ws.Range("A12").Value = "Breached"
$ws.Range("J12").Value = "~"

# --- Column J no longer needs its own custom width now that its content
# collapsed back to the default (mirrors the source width bookkeeping:
# column 10 drops out of the explicit <cols> custom-width list). Clear
# the column-level formatting, then restore the original per-cell look
# (bold+centered header, centered body) so only the width reverts.
$ws.Columns.Item(10).ClearFormats()
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").Font.Bold = $true
$ws.Range("J2:J13").HorizontalAlignment = -4108
